# Trade #42 closed at 2026-02-17 21:04:07 - unknown UNKNOWN +0.000%
#
# This script applies the following changes to the live trading results
# workbook:
#   1. Summary sheet: Total Trades 69 -> 70, Win Rate % 46.38 -> 45.71
#   2. Strategy Status sheet: MarketMaking row - Trades 36 -> 37,
#      Win Rate % 50 -> 48.65
#   3. All Trades sheet: close out existing open trade (row 71, trade #70)
#      as CLOSED via early_exit, and append a brand-new open trade
#      (row 104, trade #103)
#   4. MarketMaking sheet: same two edits mirrored in its own layout
#      (row 38, trade #70 closed; row 71, trade #103 appended)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 70
$summary.Range("B9").Value = 45.71

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking is row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 37
$status.Range("G5").Value = 48.65

# ---------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out the existing open trade in row 71 (Trade # 70)
$allTrades.Range("G71").Value = 0.98
$allTrades.Range("H71").Value = "CLOSED"
$allTrades.Range("K71").Value = 100.56
$allTrades.Range("L71").Value = "early_exit"
$allTrades.Range("M71").Value = 0.11

# Append the newly opened trade as row 104 (Trade # 103)
$allTrades.Range("A104").Value = 103
$allTrades.Range("B104").NumberFormat = "@"
$allTrades.Range("B104").Value = "2026-02-17"
$allTrades.Range("C104").NumberFormat = "@"
$allTrades.Range("C104").Value = "21:04:00"
$allTrades.Range("D104").Value = "MarketMaking"
$allTrades.Range("E104").Value = "DOWN"
$allTrades.Range("F104").Value = 0.98
$allTrades.Range("H104").Value = "OPEN"
$allTrades.Range("I104").Value = 0
$allTrades.Range("J104").Value = 0
$allTrades.Range("K104").Value = 100.5619219857093
$allTrades.Range("M104").Value = 0
$allTrades.Range("N104").Value = 0
$allTrades.Range("O104").Value = 0
$allTrades.Range("P104").Value = 0.6
$allTrades.Range("Q104").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# 4. MarketMaking sheet (per-strategy log, different column order)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close out the existing open trade in row 38 (Trade # 70)
$mm.Range("G38").Value = 0.98
$mm.Range("H38").Value = "CLOSED"
$mm.Range("K38").Value = 100.56
$mm.Range("P38").Value = "early_exit"
$mm.Range("Q38").Value = 0.11

# Append the newly opened trade as row 71 (Trade # 103)
$mm.Range("A71").Value = 103
$mm.Range("B71").NumberFormat = "@"
$mm.Range("B71").Value = "2026-02-17"
$mm.Range("C71").NumberFormat = "@"
$mm.Range("C71").Value = "21:04:00"
$mm.Range("D71").Value = "MarketMaking"
$mm.Range("E71").Value = "DOWN"
$mm.Range("F71").Value = 0.98
$mm.Range("H71").Value = "OPEN"
$mm.Range("I71").Value = 0
$mm.Range("J71").Value = 0
$mm.Range("K71").Value = 100.5619219857093
$mm.Range("L71").Value = 0
$mm.Range("M71").Value = 0
$mm.Range("N71").Value = 0.6
$mm.Range("O71").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q71").Value = 0
